$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices in column D that parse as plain numbers (e.g. "1.002", "335.85")
# must keep their original text representation (incl. trailing zeros),
# so force those specific cells to Text format before assigning the value.
# (Values like "27.176.86" already have two dots and can never be
# misread as a number, so they are left at their default format.)
$ws.Range("D4:D14").NumberFormat = "@"
$ws.Range("D16:D22").NumberFormat = "@"
$ws.Range("D24:D29").NumberFormat = "@"
$ws.Range("D31:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '27.176.86'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.780.85'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '335.85'
$ws.Range("E5").Value = '  -1.85%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.3774'
$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3425'
$ws.Range("E8").Value = '  -2.42%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '48.20'
$ws.Range("E9").Value = '  -2.73%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.192'
$ws.Range("E10").Value = '  -2.43%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.07436'
$ws.Range("E11").Value = '  -3.62%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '21.53'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '6.444'
$ws.Range("E14").Value = '  -2.91%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.776.73'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '7.046'
$ws.Range("E16").Value = '  -2.19%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001095'
$ws.Range("E17").Value = '  -2.21%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.06642'
$ws.Range("E18").Value = '  -2.28%  '

$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '84.44'
$ws.Range("E19").Value = '  -1.87%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '0.9989'
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.656'
$ws.Range("E21").Value = '  +3.32%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '17.23'
$ws.Range("E22").Value = '  -3.22%  '

$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '27.156.07'
$ws.Range("E23").Value = '  -1.02%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '12.35'
$ws.Range("E24").Value = '  -5.49%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.426'
$ws.Range("E25").Value = '  -1.71%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '1.533'
$ws.Range("E26").Value = '  +2.65%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.541'
$ws.Range("E27").Value = '  -1.52%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '21.42'
$ws.Range("E28").Value = '  +2.53%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '151.66'
$ws.Range("E29").Value = '  -1.59%  '

$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '1.981.31'
$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '133.30'
$ws.Range("E31").Value = '  -2.36%  '

$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '4.062'
$ws.Range("E32").Value = '  -3.01%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.026'
$ws.Range("E33").Value = '  -4.87%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '0.08617'
$ws.Range("E34").Value = '  -2.14%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '13.09'
$ws.Range("E35").Value = '  -5.42%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.657'
$ws.Range("E36").Value = '  -3.42%  '

$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").Value = '0.6885'
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.410'
$ws.Range("E38").Value = '  -4.42%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06335'
$ws.Range("E39").Value = '  -3.41%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.746'
$ws.Range("E40").Value = '  -2.62%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2181'
$ws.Range("E41").Value = '  -3.62%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.02321'
$ws.Range("E42").Value = '  -4.45%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.239'
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '14.43'
$ws.Range("E44").Value = '  -1.99%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '0.9996'
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6397'
$ws.Range("E46").Value = '  -0.86%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.842'
$ws.Range("E47").Value = '  -4.40%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.129'
$ws.Range("E48").Value = '  -2.76%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '128.76'
$ws.Range("E49").Value = '  -3.07%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.07113'
$ws.Range("E50").Value = '  -3.35%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '79.06'
$ws.Range("E51").Value = '  -2.17%  '
